$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) CC-EUpTIECS ("CC-EUpTIECS"): add "green hydrogen" / "low carbon hydrogen"
#    columns (L, M) to the Energy Use per Ton Sequestered table.
#    NOTE: this sheet is filled out first so the two new shared strings get
#    allocated shared-string indexes 98/99 (ahead of the "...if" variants
#    added below on CC-EUpTHSCS).
# ---------------------------------------------------------------------------
$wsIE = $wb.Worksheets.Item("CC-EUpTIECS")

$wsIE.Range("L1").Value = "green hydrogen"
$wsIE.Range("M1").Value = "low carbon hydrogen"
$wsIE.Range("K1").Copy()
$wsIE.Range("L1:M1").PasteSpecial(-4122)   # xlPasteFormats

for ($r = 2; $r -le 26; $r++) {
    $wsIE.Cells.Item($r, 12).Value = 0
    $wsIE.Cells.Item($r, 13).Value = 0
}

# ---------------------------------------------------------------------------
# 2) CC-EUpTIPCS: same green/low carbon hydrogen columns (L, M).
# ---------------------------------------------------------------------------
$wsIP = $wb.Worksheets.Item("CC-EUpTIPCS")

$wsIP.Range("L1").Value = "green hydrogen"
$wsIP.Range("M1").Value = "low carbon hydrogen"
$wsIP.Range("K1").Copy()
$wsIP.Range("L1:M1").PasteSpecial(-4122)   # xlPasteFormats

for ($r = 2; $r -le 26; $r++) {
    $wsIP.Cells.Item($r, 12).Value = 0
    $wsIP.Cells.Item($r, 13).Value = 0
}

# Selection on CC-EUpTIPCS moves to the freshly added L:M columns.
$wsIP.Range("L1:M1048576").Select()

# ---------------------------------------------------------------------------
# 3) CC-EUpTHSCS: add "green hydrogen if" / "low carbon hydrogen if" rows
#    (12, 13) feeding the Hydrogen Sector demand-if subscript list.
# ---------------------------------------------------------------------------
$wsHS = $wb.Worksheets.Item("CC-EUpTHSCS")

$wsHS.Range("A12").Value = "green hydrogen if"
$wsHS.Range("B12").Value = 0
$wsHS.Range("A13").Value = "low carbon hydrogen if"
$wsHS.Range("B13").Value = 0

$wsHS.Range("A11").Copy()
$wsHS.Range("A12:A13").PasteSpecial(-4122) # xlPasteFormats

# ---------------------------------------------------------------------------
# 4) CC-CCoIPCE: selection only moves from B3 to C3 (no data change).
# ---------------------------------------------------------------------------
$wsCP = $wb.Worksheets.Item("CC-CCoIPCE")
$wsCP.Range("C3").Select()

# ---------------------------------------------------------------------------
# 5) CC-EUpTHSCS becomes the active sheet / tab in the workbook. Activating
#    (and re-selecting its range) LAST ensures it stays the active tab even
#    though other sheets' ranges were selected above.
# ---------------------------------------------------------------------------
$wsHS.Activate()
$wsHS.Range("A14").Select()
